# Insert a new data row at row 55 (pushing the existing rows 55-150 down to
# 56-151) and populate it with the new daily-price record, as described by
# the commit "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 55..150 down to 56..151, leaving a blank row 55.
$ws.Rows(55).Insert()

# Populate the newly inserted row 55 with the new record.
$ws.Range("A55").Value = 4
$ws.Range("B55").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C55").Value = "Los Lagos"
$ws.Range("D55").Value = 44469
$ws.Range("E55").Value = 10
$ws.Range("F55").Value = 100112043
$ws.Range("G55").Value = "Pepino ensalada"
$ws.Range("H55").Value = "Sin especificar"
$ws.Range("I55").Value = "Primera"
$ws.Range("J55").Value = 200
$ws.Range("K55").Value = 20000
$ws.Range("L55").Value = 20000
$ws.Range("M55").Value = 20000
$ws.Range("N55").Value = "$/caja 60 unidades"
$ws.Range("O55").Value = "Región de Arica y Parinacota"
$ws.Range("P55").Value = 333
$ws.Range("Q55").Value = 60
$ws.Range("R55").Value = "Hortaliza"
